# "Add chinese html tag"
# The "updates" sheet had two extra duplicate rows (the hello-中/CVEnumber
# block was repeated 3x) that get trimmed down to a single row, and every
# ":" used as a separator inside CVE/version strings is swapped for "@".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("updates")

# Row 5 keeps the "hello中" / "CVEnumber..." data but its risk value (F5)
# moves from "高" to what used to live in F7 ("3").
$ws.Range("F5").Value = "3"

# Drop the two now-redundant duplicate rows (old rows 6 and 7); row 8
# shifts up into row 6.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# Swap ":" for "@" in the CVE / version columns (column B) across all
# the data rows that used that separator.
$ws.Range("B3").Value = "b@1"
$ws.Range("B5").Value = "CVEnumber@CVEurl,2@asdf,3@123123"
$ws.Range("B2").Value = "02@222"
$ws.Range("B4").Value = "b@1"
$ws.Range("B6").Value = "b@1"

# Match the saved selection state (B11 -> B6, since the sheet shrank).
$ws.Range("B6").Select()
